$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at the very start of the sheet. Everything that
# used to live in A:E (Class, Method/Prop/Rel Name, Member ID, SnippetId,
# MethodName) shifts one column to the right, landing in B:F.
$ws.Columns.Item(1).Insert()

# Give the new first column a header and fill its data rows with the
# package name every existing row belongs to.
$ws.Range("A1").Value = "Package"
$ws.Range("A2:A12").Value = "PowerPoint"

# The existing "Snippets" table still only covers the old B:F extent (and
# its column metadata doesn't auto-shift with the raw column insert), so
# drop it and recreate it over the full A1:F12 range. Excel will read the
# column names straight back out of row 1, giving us Package/Class/
# Method-Prop-Rel Name/Member ID/SnippetId/MethodName in the right order.
$lo = $ws.ListObjects.Item(1)
$lo.Unlist()

$lo2 = $ws.ListObjects.Add(1, $ws.Range("A1:F12"), [System.Reflection.Missing]::Value, 1)
$lo2.Name = "Snippets"
$lo2.TableStyle = "TableStyleMedium7"

# The data rows (3-12) carried an explicit style index purely so Excel
# could track "no special formatting" for the table body; clear it back to
# the workbook default now that the table has been rebuilt.
$ws.Range("B3:F12").Style = "Normal"

# Leave the same cells selected that were just populated.
$ws.Range("A2:A12").Select() | Out-Null
